$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple single-run text edits (rows 1-4) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "1395"

# --- Row 6: 0.00021 -> 0.00091, then insert a brand-new row (0.00024) right after it ---
$row6 = $t.Rows.Item(6)
$row6.Cells.Item(1).Range.Text = "0.00091"
$newRow = $t.Rows.Add($t.Rows.Item(7))
$newRow.Cells.Item(1).Range.Text = "0.00024"

# --- Edits to the rows that followed the old row 6 (now shifted down by one
#     because of the insert above); grab them by object reference so further
#     structural edits (the coming delete) can't desync the indices. ---
$rowA = $t.Rows.Item(9)   # was old row 8: 0.00002 -> 0.00039
$rowA.Cells.Item(1).Range.Text = "0.00039"

$rowB = $t.Rows.Item(10)  # was old row 9: 0.00006 -> 0.00051
$rowB.Cells.Item(1).Range.Text = "0.00051"

$rowC = $t.Rows.Item(11)  # was old row 10: 0.00007 -> 0.00068
$rowC.Cells.Item(1).Range.Text = "0.00068"

# --- Old row 11 (0.00009) is removed entirely; old row 12 (0.02472) becomes 0.34072 ---
$rowToDelete = $t.Rows.Item(12)   # old row 11 (0.00009), now at index 12
$rowToDelete.Delete()

$rowRenamed = $t.Rows.Item(12)    # old row 12 (0.02472), now settled at index 12
$rowRenamed.Cells.Item(1).Range.Text = "0.34072"

# --- Collapse the three multi-run rows near the end of the table into single runs ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.93"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.34"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "497"

$d.Save()
